# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Betarraga"
# at row 251, pushing the existing rows 251-318 down to 252-319.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 251 (shifts rows 251..318 down to 252..319)
$ws.Rows("251:251").Insert()

# Populate the newly inserted row with the new record
$ws.Range('A251').Value = 4
$ws.Range('B251').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C251').Value = 'Los Lagos'
$ws.Range('D251').Value = 44754
$ws.Range('E251').Value = 10
$ws.Range('F251').Value = 100114014
$ws.Range('G251').Value = 'Betarraga'
$ws.Range('H251').Value = 'Sin especificar'
$ws.Range('I251').Value = 'Primera'
$ws.Range('J251').Value = 800
$ws.Range('K251').Value = 1200
$ws.Range('L251').Value = 1200
$ws.Range('M251').Value = 1200
$ws.Range('N251').Value = '$/paquete 5 unidades'
$ws.Range('O251').Value = 'Región del Maule'
$ws.Range('P251').Value = 240
$ws.Range('Q251').Value = 5
$ws.Range('R251').Value = 'Hortaliza'
